# Append updated data rows (302-328) to Sheet1, mirroring the formatting
# of the last existing row (301), and update the sheet's used dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data: row, dateSerial(colA), colB, colC, colD
$data = @(
    @(302,44376,0,0,0),
    @(303,44377,1,1,6.628222973420826),
    @(304,44378,1,2,13.25644594684165),
    @(305,44379,0,2,13.25644594684165),
    @(306,44380,0,2,13.25644594684165),
    @(307,44381,0,2,13.25644594684165),
    @(308,44382,0,2,13.25644594684165),
    @(309,44383,0,2,13.25644594684165),
    @(310,44384,1,2,13.25644594684165),
    @(311,44385,0,1,6.628222973420826),
    @(312,44386,0,1,6.628222973420826),
    @(313,44387,0,1,6.628222973420826),
    @(314,44388,0,1,6.628222973420826),
    @(315,44389,1,2,13.25644594684165),
    @(316,44390,0,2,13.25644594684165),
    @(317,44391,0,1,6.628222973420826),
    @(318,44392,0,1,6.628222973420826),
    @(319,44393,0,1,6.628222973420826),
    @(320,44394,0,1,6.628222973420826),
    @(321,44395,1,2,13.25644594684165),
    @(322,44396,0,1,6.628222973420826),
    @(323,44397,0,1,6.628222973420826),
    @(324,44398,0,1,6.628222973420826),
    @(325,44399,2,3,19.88466892026248),
    @(326,44400,0,3,19.88466892026248),
    @(327,44401,0,3,19.88466892026248),
    @(328,44402,1,3,19.88466892026248)
)

$lastRow = 301

foreach ($entry in $data) {
    $r = $entry[0]
    $dateSerial = $entry[1]
    $colB = $entry[2]
    $colC = $entry[3]
    $colD = $entry[4]

    # Copy the formatting (number format, style, borders, alignment) from the
    # corresponding cell in the last existing row down to the new row.
    $ws.Range("A$lastRow").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("B$lastRow`:D$lastRow").Copy() | Out-Null
    $ws.Range("B$r`:D$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 4).Value = $colD
}

$excel.CutCopyMode = 0
